# This script applies the cryptocurrency price/volume update described by the commit
# "Updated cryptos list on Fri Nov 17 10:13:43 UTC 2023 with GitHub Actions".
#
# For each coin row, the Price (column D) and Volume(1h) (column E) values are refreshed.
# Two coin pairs (rows 22/23 and rows 37/38) swapped rank order, so their Coin name (B),
# Link (C), Price (D) and Volume (E) values are exchanged together with the new figures.
#
# Numeric-looking Price strings (e.g. "5.17") are written with a leading apostrophe so
# Excel stores them as text (matching the original inline-string/text cell type) instead
# of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.290.99"
$ws.Range("E2").Value = "  -2.99%  "
$ws.Range("D3").Value = "1.968.86"
$ws.Range("E3").Value = "  -3.71%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'244.18"
$ws.Range("E5").Value = "  -3.19%  "
$ws.Range("D6").Value = "'0.623"
$ws.Range("E6").Value = "  -4.34%  "
$ws.Range("D7").Value = "'57.40"
$ws.Range("E7").Value = "  -13.38%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.371"
$ws.Range("E9").Value = "  -7.13%  "
$ws.Range("D10").Value = "'55.86"
$ws.Range("E10").Value = "  -5.67%  "
$ws.Range("D11").Value = "'0.0867"
$ws.Range("E11").Value = "  +8.43%  "
$ws.Range("E12").Value = "  -0.20%  "
$ws.Range("D13").Value = "'0.840"
$ws.Range("E13").Value = "  -8.23%  "
$ws.Range("D14").Value = "'21.79"
$ws.Range("E14").Value = "  -7.51%  "
$ws.Range("D15").Value = "2.261.24"
$ws.Range("E15").Value = "  -3.56%  "
$ws.Range("D16").Value = "'13.58"
$ws.Range("E16").Value = "  -8.49%  "
$ws.Range("D17").Value = "'5.36"
$ws.Range("E17").Value = "  -6.50%  "
$ws.Range("D18").Value = "1.962.27"
$ws.Range("E18").Value = "  -3.91%  "
$ws.Range("D19").Value = "36.187.27"
$ws.Range("E19").Value = "  -2.93%  "
$ws.Range("D20").Value = "0.0₃0893"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("D21").Value = "'70.88"
$ws.Range("E21").Value = "  -3.12%  "
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'233.12"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'5.17"
$ws.Range("E23").Value = "  -6.60%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -3.74%  "
$ws.Range("D26").Value = "'2.28"
$ws.Range("E26").Value = "  -4.14%  "
$ws.Range("D27").Value = "'9.66"
$ws.Range("E27").Value = "  -3.45%  "
$ws.Range("D28").Value = "'165.78"
$ws.Range("E28").Value = "  +2.43%  "
$ws.Range("D29").Value = "'20.21"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").Value = "'0.127"
$ws.Range("E30").Value = "  -2.06%  "
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("D32").Value = "'1.16"
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").Value = "'4.79"
$ws.Range("E33").Value = "  -7.04%  "
$ws.Range("D34").Value = "'0.0644"
$ws.Range("E34").Value = "  +2.00%  "
$ws.Range("D35").Value = "'4.36"
$ws.Range("E35").Value = "  -6.45%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.80"
$ws.Range("E37").Value = "  -2.08%  "
$ws.Range("B38").Value = "THORChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D38").Value = "'5.99"
$ws.Range("E38").Value = "  -6.61%  "
$ws.Range("D39").Value = "'2.16"
$ws.Range("E39").Value = "  -8.84%  "
$ws.Range("D40").Value = "'2.91"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("D41").Value = "'0.0958"
$ws.Range("E41").Value = "  -5.60%  "
$ws.Range("E42").Value = "  -7.85%  "
$ws.Range("D43").Value = "'2.88"
$ws.Range("E43").Value = "  -5.18%  "
$ws.Range("D44").Value = "'0.0211"
$ws.Range("E44").Value = "  -3.63%  "
$ws.Range("D45").Value = "'1.07"
$ws.Range("E45").Value = "  -8.33%  "
$ws.Range("D46").Value = "'15.83"
$ws.Range("E46").Value = "  -9.64%  "
$ws.Range("D47").Value = "'89.34"
$ws.Range("E47").Value = "  -6.47%  "
$ws.Range("D48").Value = "1.350.84"
$ws.Range("E48").Value = "  -3.20%  "
$ws.Range("D49").Value = "'7.32"
$ws.Range("E49").Value = "  -6.53%  "
$ws.Range("D50").Value = "'2.81"
$ws.Range("E50").Value = "  -3.57%  "
$ws.Range("D51").Value = "'44.66"
$ws.Range("E51").Value = "  -5.42%  "
